# heatmap_regions.xlsx maintenance edit:
#   - rename the "Mobotrex" sheet tab to "MoboTrex"
#   - reword the install/quantity comment strings to mention "ATC installations"
#   - refresh the stale active-cell selections left over on a few tabs
#   - nudge the iterative-calculation max-change setting and the saved
#     window geometry to match the author's last session (best effort —
#     see notes at bottom)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet rename: "Mobotrex" -> "MoboTrex"
# ---------------------------------------------------------------------
$wsMoboTrex = $wb.Worksheets.Item("Mobotrex")
$wsMoboTrex.Name = "MoboTrex"

# ---------------------------------------------------------------------
# 2. Reworded comment text ("... installed" -> "ATC installations ...")
# ---------------------------------------------------------------------
$wsEconolite = $wb.Worksheets.Item("Econolite")
$wsEconolite.Range("D3").Value = "ATC installations in 200+ intersections - 10/30/2023"

$wsQFree = $wb.Worksheets.Item("Q-Free")
$wsQFree.Range("D2").Value = "269 ATC installations - 04/08/21"
$wsQFree.Range("D3").Value = "800+ ATC installations - 12/12/2023"

$wsCubic = $wb.Worksheets.Item("Cubic")
$wsCubic.Range("D2").Value = "ATC installations in ~20 intersections - 10/03/17"

# ---------------------------------------------------------------------
# 3. Per-sheet selection / active-cell updates
#    (activating each sheet in turn records its own selection; the last
#    sheet activated becomes the workbook's active tab, so we finish on
#    "Western Systems" to match the saved tabSelected/activeTab state)
# ---------------------------------------------------------------------
[void]$wsEconolite.Activate()
[void]$wsEconolite.Range("D13").Select()

[void]$wsQFree.Activate()
[void]$wsQFree.Range("E18").Select()

[void]$wsCubic.Activate()
[void]$wsCubic.Range("D2").Select()

$wsWestern = $wb.Worksheets.Item("Western Systems")
[void]$wsWestern.Activate()
[void]$wsWestern.Range("H27").Select()

# ---------------------------------------------------------------------
# 4. Iterative-calculation max-change (workbook calcPr iterateDelta).
#    Excel exposes this as Application.Iteration / MaxIterations /
#    MaxChange; set them to mirror the authored value (1E-4) even though
#    this headless host's xlsx writer does not currently round-trip the
#    calcPr@iterateDelta attribute.
# ---------------------------------------------------------------------
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

# ---------------------------------------------------------------------
# 5. Saved window geometry (bookViews/workbookView yWindow & windowHeight).
#    Mirrors the author's last-saved window position/size; this host does
#    not persist Window.Top/Height back into workbookView, but setting
#    them keeps intent explicit and is harmless if unsupported.
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Top = 500
$win.Height = 21100
